$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap rows 8 <-> 9 (A, Q, R columns; B/D/E/F/G/H already identical)
$r8A = $ws.Range("A8").Value2
$r8Q = $ws.Range("Q8").Value2
$r8R = $ws.Range("R8").Value2
$r9A = $ws.Range("A9").Value2
$r9Q = $ws.Range("Q9").Value2
$r9R = $ws.Range("R9").Value2

$ws.Range("A8").Value2 = $r9A
$ws.Range("Q8").Value2 = $r9Q
$ws.Range("R8").Value2 = $r9R
$ws.Range("A9").Value2 = $r8A
$ws.Range("Q9").Value2 = $r8Q
$ws.Range("R9").Value2 = $r8R

# Swap rows 10 <-> 11 (A, B, E, F, G, H, Q, R columns; D identical)
$r10A = $ws.Range("A10").Value2
$r10B = $ws.Range("B10").Value2
$r10E = $ws.Range("E10").Value2
$r10F = $ws.Range("F10").Value2
$r10G = $ws.Range("G10").Value2
$r10H = $ws.Range("H10").Value2
$r10Q = $ws.Range("Q10").Value2
$r10R = $ws.Range("R10").Value2

$r11A = $ws.Range("A11").Value2
$r11B = $ws.Range("B11").Value2
$r11E = $ws.Range("E11").Value2
$r11F = $ws.Range("F11").Value2
$r11G = $ws.Range("G11").Value2
$r11H = $ws.Range("H11").Value2
$r11Q = $ws.Range("Q11").Value2
$r11R = $ws.Range("R11").Value2

$ws.Range("A10").Value2 = $r11A
$ws.Range("B10").Value2 = $r11B
$ws.Range("E10").Value2 = $r11E
$ws.Range("F10").Value2 = $r11F
$ws.Range("G10").Value2 = $r11G
$ws.Range("H10").Value2 = $r11H
$ws.Range("Q10").Value2 = $r11Q
$ws.Range("R10").Value2 = $r11R

$ws.Range("A11").Value2 = $r10A
$ws.Range("B11").Value2 = $r10B
$ws.Range("E11").Value2 = $r10E
$ws.Range("F11").Value2 = $r10F
$ws.Range("G11").Value2 = $r10G
$ws.Range("H11").Value2 = $r10H
$ws.Range("Q11").Value2 = $r10Q
$ws.Range("R11").Value2 = $r10R

# Swap rows 13 <-> 14 (A, B, D, E, F, G, H, Q, R columns)
$r13A = $ws.Range("A13").Value2
$r13B = $ws.Range("B13").Value2
$r13D = $ws.Range("D13").Value2
$r13E = $ws.Range("E13").Value2
$r13F = $ws.Range("F13").Value2
$r13G = $ws.Range("G13").Value2
$r13H = $ws.Range("H13").Value2
$r13Q = $ws.Range("Q13").Value2
$r13R = $ws.Range("R13").Value2

$r14A = $ws.Range("A14").Value2
$r14B = $ws.Range("B14").Value2
$r14D = $ws.Range("D14").Value2
$r14E = $ws.Range("E14").Value2
$r14F = $ws.Range("F14").Value2
$r14G = $ws.Range("G14").Value2
$r14H = $ws.Range("H14").Value2
$r14Q = $ws.Range("Q14").Value2
$r14R = $ws.Range("R14").Value2

$ws.Range("A13").Value2 = $r14A
$ws.Range("B13").Value2 = $r14B
$ws.Range("D13").Value2 = $r14D
$ws.Range("E13").Value2 = $r14E
$ws.Range("F13").Value2 = $r14F
$ws.Range("G13").Value2 = $r14G
$ws.Range("H13").Value2 = $r14H
$ws.Range("Q13").Value2 = $r14Q
$ws.Range("R13").Value2 = $r14R

$ws.Range("A14").Value2 = $r13A
$ws.Range("B14").Value2 = $r13B
$ws.Range("D14").Value2 = $r13D
$ws.Range("E14").Value2 = $r13E
$ws.Range("F14").Value2 = $r13F
$ws.Range("G14").Value2 = $r13G
$ws.Range("H14").Value2 = $r13H
$ws.Range("Q14").Value2 = $r13Q
$ws.Range("R14").Value2 = $r13R
